# Update the first question in the testing matrix to reflect the new
# functionality: writing user-supplied step info to a csv/txt file at runtime.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Programul poate scrie la rulare informatiile pe care user-ul i le ofera intr-un fisier csv sau txt?"
